$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above row 5 (shifts old rows 5-9 down to 6-10).
$ws.Rows("5").Insert()

# The plain Insert() leaves the new row with no formatting/height, so copy
# row 4's formatting down into the new row 5 (matches Excel's default
# "Insert Copied Cells" / Format-From-Above behaviour the author relied on).
$ws.Range("B4:G4").Copy()
$ws.Range("B5:G5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K4:O4").Copy()
$ws.Range("K5:O5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Rows("5").RowHeight = 28.8

# New row 5 (B:G): "search for properties" sub-row under the "properties" resource
$ws.Range("B5").Value = "properties"
$ws.Range("G5").Value = "search for properties"

# New row 5 (K:O) inherits what used to be row 4's scenario
# ("/properties/[propID]" -> "Retrieve property for propID")
$ws.Range("K5").Value = "/properties/[propID]"
$ws.Range("M5").Value = "Retrieve property for propID"

# Row 4 (K:O) becomes the brand new "/properties" search scenario
$ws.Range("K4").Value = "/properties"
$ws.Range("M4").Value = "Retrieve properties filtered by params"

# Update the active selection to match the authored state
$ws.Range("H6").Select()
